$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44229
$ws.Range("K4").Value = 'Artic Sprite'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("Q4").Value = '$/bandeja 18 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1083

# Row 5
$ws.Range("D5").Value = 44223
$ws.Range("K5").Value = 'Ruby Diamond'
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1139

# Row 6
$ws.Range("D6").Value = 44223
$ws.Range("K6").Value = 'Super Queen'
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 1139

# Row 7
$ws.Range("D7").Value = 44244
$ws.Range("K7").Value = 'Nectar Crest'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19500
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1083

# Row 8
$ws.Range("D8").Value = 44244
$ws.Range("K8").Value = 'Venus'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1083

# Row 9
$ws.Range("D9").Value = 44202
$ws.Range("K9").Value = 'Super Queen'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1167

# Row 10
$ws.Range("D10").Value = 44523
$ws.Range("K10").Value = 'Early Glo'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 24000
$ws.Range("O10").Value = 25000
$ws.Range("P10").Value = 24500
$ws.Range("Q10").Value = '$/bandeja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1361

# Row 11
$ws.Range("D11").Value = 44524
$ws.Range("K11").Value = 'Early Glo'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 24500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("R11").Value = 'Región de Coquimbo'
$ws.Range("S11").Value = 1361

# Row 12
$ws.Range("D12").Value = 44216
$ws.Range("K12").Value = 'Nectar Crest'
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19500
$ws.Range("Q12").Value = '$/bandeja 18 kilos granel'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 1083

# Row 13
$ws.Range("D13").Value = 44169
$ws.Range("K13").Value = 'Artic Sprite'
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 24500
$ws.Range("Q13").Value = '$/bandeja 18 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 1361

# Row 14
$ws.Range("D14").Value = 44169
$ws.Range("K14").Value = 'Early John'
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 270
$ws.Range("N14").Value = 24000
$ws.Range("O14").Value = 25000
$ws.Range("P14").Value = 24500
$ws.Range("Q14").Value = '$/bandeja 18 kilos granel'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 1361

# Row 15
$ws.Range("D15").Value = 44215
$ws.Range("K15").Value = 'Venus'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 19000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 19500
$ws.Range("Q15").Value = '$/bandeja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 1083

# Row 16
$ws.Range("D16").Value = 44201
$ws.Range("K16").Value = 'Super Queen'
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 22500
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 1250

# Row 17
$ws.Range("D17").Value = 44236
$ws.Range("K17").Value = 'June Pearl'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 270
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 21000
$ws.Range("P17").Value = 20500
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1139

# Row 18
$ws.Range("D18").Value = 44222
$ws.Range("K18").Value = 'Nectar Crest'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 270
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 21000
$ws.Range("P18").Value = 20500
$ws.Range("Q18").Value = '$/bandeja 18 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1139

# Row 19
$ws.Range("D19").Value = 44243
$ws.Range("K19").Value = 'Venus'
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 250
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 21000
$ws.Range("P19").Value = 20500
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 1139

# Row 20
$ws.Range("D20").Value = 44174
$ws.Range("K20").Value = 'Early John'
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 20500
$ws.Range("Q20").Value = '$/caja 18 kilos granel'
$ws.Range("R20").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S20").Value = 1139

# Row 21
$ws.Range("D21").Value = 44273
$ws.Range("K21").Value = 'Artic Snow'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 250
$ws.Range("N21").Value = 22000
$ws.Range("O21").Value = 23000
$ws.Range("P21").Value = 22500
$ws.Range("Q21").Value = '$/bandeja 18 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 1250

# Row 22
$ws.Range("D22").Value = 44273
$ws.Range("K22").Value = 'August Red'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 23000
$ws.Range("P22").Value = 22500
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1250

# Row 23
$ws.Range("D23").Value = 44209
$ws.Range("K23").Value = 'Super Queen'
$ws.Range("L23").Value = 'Tercera'
$ws.Range("M23").Value = 320
$ws.Range("N23").Value = 17000
$ws.Range("O23").Value = 18000
$ws.Range("P23").Value = 17500
$ws.Range("Q23").Value = '$/caja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 972

# Row 24
$ws.Range("D24").Value = 44533
$ws.Range("K24").Value = 'Artic Pride'
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 270
$ws.Range("N24").Value = 24000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 24500
$ws.Range("Q24").Value = '$/bandeja 18 kilos granel'
$ws.Range("R24").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S24").Value = 1361

# Row 30
$ws.Range("D30").Value = 44167
$ws.Range("K30").Value = 'Early John'
$ws.Range("L30").Value = 'Segunda'
$ws.Range("M30").Value = 200
$ws.Range("N30").Value = 25000
$ws.Range("O30").Value = 26000
$ws.Range("P30").Value = 25500
$ws.Range("Q30").Value = '$/caja 18 kilos granel'
$ws.Range("R30").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S30").Value = 1417

# Row 31
$ws.Range("D31").Value = 44161
$ws.Range("K31").Value = 'Artic Glo'
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 280
$ws.Range("N31").Value = 25000
$ws.Range("O31").Value = 26000
$ws.Range("P31").Value = 25500
$ws.Range("Q31").Value = '$/bandeja 18 kilos granel'
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 1417

# Row 32
$ws.Range("D32").Value = 44161
$ws.Range("K32").Value = 'Early John'
$ws.Range("L32").Value = 'Segunda'
$ws.Range("M32").Value = 250
$ws.Range("N32").Value = 26000
$ws.Range("O32").Value = 27000
$ws.Range("P32").Value = 26500
$ws.Range("Q32").Value = '$/caja 18 kilos granel'
$ws.Range("R32").Value = 'Región de O''Higgins'
$ws.Range("S32").Value = 1472

# Row 35
$ws.Range("D35").Value = 44238
$ws.Range("K35").Value = 'August Red'
$ws.Range("L35").Value = 'Segunda'
$ws.Range("M35").Value = 320
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 21000
$ws.Range("P35").Value = 20500
$ws.Range("Q35").Value = '$/caja 18 kilos granel'
$ws.Range("R35").Value = 'Región de O''Higgins'
$ws.Range("S35").Value = 1139

# Row 36
$ws.Range("D36").Value = 44238
$ws.Range("K36").Value = 'Venus'
$ws.Range("L36").Value = 'Segunda'
$ws.Range("M36").Value = 320
$ws.Range("N36").Value = 20000
$ws.Range("O36").Value = 21000
$ws.Range("P36").Value = 20500
$ws.Range("Q36").Value = '$/bandeja 18 kilos granel'
$ws.Range("R36").Value = 'Región de O''Higgins'
$ws.Range("S36").Value = 1139

# Row 37
$ws.Range("D37").Value = 44257
$ws.Range("K37").Value = 'August Red'
$ws.Range("L37").Value = 'Segunda'
$ws.Range("M37").Value = 300
$ws.Range("N37").Value = 19000
$ws.Range("O37").Value = 20000
$ws.Range("P37").Value = 19500
$ws.Range("Q37").Value = '$/caja 18 kilos granel'
$ws.Range("R37").Value = 'Región de O''Higgins'
$ws.Range("S37").Value = 1083

# Row 38
$ws.Range("D38").Value = 44540
$ws.Range("K38").Value = 'Artic Pride'
$ws.Range("L38").Value = 'Segunda'
$ws.Range("M38").Value = 250
$ws.Range("N38").Value = 21000
$ws.Range("O38").Value = 22000
$ws.Range("P38").Value = 21500
$ws.Range("Q38").Value = '$/bandeja 18 kilos granel'
$ws.Range("R38").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S38").Value = 1194
